# "chasing qpsk TX distort": retune the "Branch LPF" sheet's cutoff/sample
# rate inputs (B2, B3). Everything else on the sheet (B4..B42) is derived
# via formulas, so Excel's automatic recalculation produces the rest of the
# diff once these two inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Branch LPF")

$ws.Range("B2").Value = 150
$ws.Range("B3").Value = 86400

# Make the sheet active and move the selection to B3 (matches the author's
# resulting view state: scrolled back to the top with B3 selected instead
# of the previously scrolled-down D43 selection).
$ws.Activate()
$ws.Range("B3").Select()

$wb.Save()
